$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "A"
$ws.Range("B6").Value = "A"
$ws.Range("B7").Value = "A"
$ws.Range("B8").Value = "BB"
$ws.Range("B9").Value = "AAA"
$ws.Range("B10").Value = "AAA"
$ws.Range("B13").Value = "BBB"
$ws.Range("B14").Value = "BBB"
$ws.Range("B15").Value = "BBB"
$ws.Range("B19").Value = "A"
$ws.Range("B20").Value = "BB"
$ws.Range("B28").Value = "A"
$ws.Range("B33").Value = "BBB"
$ws.Range("B36").Value = "BBB"
$ws.Range("B37").Value = "BBB"
$ws.Range("B39").Value = "BBB"
$ws.Range("B43").Value = "BB"
$ws.Range("B45").Value = "BBB"
$ws.Range("B47").Value = "BB"
$ws.Range("B49").Value = "A"
$ws.Range("B52").Value = "B"
$ws.Range("B56").Value = "A"
$ws.Range("B59").Value = "BB"
$ws.Range("B60").Value = "BB"
$ws.Range("B61").Value = "A"
$ws.Range("B64").Value = "BBB"
$ws.Range("B65").Value = "BBB"
$ws.Range("B66").Value = "BBB"
$ws.Range("B67").Value = "BBB"
$ws.Range("B69").Value = "BB"
$ws.Range("B70").Value = "B"
$ws.Range("B72").Value = "A"
$ws.Range("B73").Value = "BBB"
$ws.Range("B75").Value = "AA"
$ws.Range("B76").Value = "BB"
$ws.Range("B77").Value = "BBB"
$ws.Range("B80").Value = "B"
$ws.Range("B82").Value = "AA"
$ws.Range("B83").Value = "AA"
$ws.Range("B86").Value = "BBB"
$ws.Range("B87").Value = "BBB"
$ws.Range("B92").Value = "BBB"
$ws.Range("B95").Value = "A"
$ws.Range("B98").Value = "A"
$ws.Range("B104").Value = "BBB"
$ws.Range("B106").Value = "AA"
$ws.Range("B107").Value = "BBB"
$ws.Range("B112").Value = "BBB"
$ws.Range("B119").Value = "BBB"
$ws.Range("B120").Value = "A"
$ws.Range("B122").Value = "BBB"
$ws.Range("B124").Value = "A"
$ws.Range("B125").Value = "A"
$ws.Range("B126").Value = "A"
$ws.Range("B129").Value = "BBB"
$ws.Range("B134").Value = "A"
$ws.Range("B141").Value = "BB"
$ws.Range("B142").Value = "A"
$ws.Range("B144").Value = "A"
$ws.Range("B146").Value = "BBB"
$ws.Range("B152").Value = "B"
$ws.Range("B154").Value = "BB"
$ws.Range("B156").Value = "B"
$ws.Range("B157").Value = "B"
$ws.Range("B158").Value = "BB"
$ws.Range("B159").Value = "AA"
$ws.Range("B160").Value = "A"
$ws.Range("B163").Value = "A"
$ws.Range("B165").Value = "B"
$ws.Range("B169").Value = "BB"
$ws.Range("B171").Value = "BB"
$ws.Range("B173").Value = "B"
$ws.Range("B178").Value = "A"
$ws.Range("B180").Value = "AA"
$ws.Range("B182").Value = "BBB"
$ws.Range("B184").Value = "BB"
$ws.Range("B185").Value = "BB"
$ws.Range("B186").Value = "AAA"
$ws.Range("B187").Value = "A"
$ws.Range("B190").Value = "BBB"
$ws.Range("B196").Value = "BB"
$ws.Range("B197").Value = "BBB"
$ws.Range("B204").Value = "A"
$ws.Range("B210").Value = "BBB"
$ws.Range("B214").Value = "BBB"
$ws.Range("B226").Value = "A"
$ws.Range("B227").Value = "BBB"
$ws.Range("B230").Value = "B"
$ws.Range("B236").Value = "BBB"
$ws.Range("B237").Value = "BB"
$ws.Range("B238").Value = "B"
$ws.Range("B240").Value = "A"
$ws.Range("B241").Value = "AAA"
$ws.Range("B252").Value = "BB"
$ws.Range("B255").Value = "BB"
$ws.Range("B256").Value = "BB"
$ws.Range("B257").Value = "BB"
$ws.Range("B258").Value = "BBB"
$ws.Range("B259").Value = "BBB"
$ws.Range("B261").Value = "A"
$ws.Range("B271").Value = "BBB"
$ws.Range("B274").Value = "A"
$ws.Range("B282").Value = "A"
$ws.Range("B284").Value = "BBB"
$ws.Range("B286").Value = "BBB"
$ws.Range("B287").Value = "A"
$ws.Range("B288").Value = "A"
$ws.Range("B290").Value = "BBB"
$ws.Range("B291").Value = "A"
$ws.Range("B293").Value = "A"
$ws.Range("B301").Value = "AA"
$ws.Range("B302").Value = "A"
$ws.Range("B303").Value = "A"
$ws.Range("B306").Value = "A"
$ws.Range("B308").Value = "A"
$ws.Range("B309").Value = "A"
$ws.Range("B314").Value = "B"
$ws.Range("B316").Value = "BBB"
$ws.Range("B318").Value = "BB"
$ws.Range("B323").Value = "B"
$ws.Range("B324").Value = "B"
$ws.Range("B326").Value = "B"
$ws.Range("B332").Value = "B"
$ws.Range("B334").Value = "B"
$ws.Range("B335").Value = "BB"
$ws.Range("B337").Value = "BB"
$ws.Range("B344").Value = "BB"
$ws.Range("B345").Value = "BB"
$ws.Range("B348").Value = "BB"
$ws.Range("B350").Value = "B"
$ws.Range("B354").Value = "BB"
$ws.Range("B356").Value = "BB"
$ws.Range("B357").Value = "BB"
$ws.Range("B358").Value = "BBB"
$ws.Range("B360").Value = "A"
$ws.Range("B361").Value = "A"
$ws.Range("B365").Value = "B"
$ws.Range("B366").Value = "A"
$ws.Range("B368").Value = "BBB"
$ws.Range("B371").Value = "BBB"
$ws.Range("B376").Value = "A"
$ws.Range("B388").Value = "B"
$ws.Range("B394").Value = "BBB"
$ws.Range("B395").Value = "BBB"
$ws.Range("B401").Value = "AAA"
$ws.Range("B402").Value = "BBB"
